$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-20 Saturday", "2025-09-21 Sunday"),
    @("35×88=", "51×86="),
    @("32×62=", "71×41="),
    @("18×79=", "31×66="),
    @("96×56=", "86×17="),
    @("45×91=", "43×23="),
    @("91×90=", "82×49="),
    @("92×97=", "50×13="),
    @("72×37=", "93×12="),
    @("57×69=", "83×64="),
    @("57×17=", "91×94="),
    @("76×20=", "45×38="),
    @("42×70=", "24×91="),
    @("35×91=", "41×59="),
    @("20×65=", "53×65="),
    @("43×78=", "66×20="),
    @("61×56=", "87×55="),
    @("44×23=", "15×34="),
    @("18×23=", "69×16="),
    @("46×88=", "39×19="),
    @("86×99=", "88×89="),
    @("66×68=", "81×58="),
    @("18×42=", "31×35="),
    @("47×42=", "84×48="),
    @("37×87=", "52×34="),
    @("42×96=", "73×62=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
